$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 16: "Borrowing or Referencing" -------------------------------------
# Copy the formatting of row 15 (plain bordered cells, style index 2) down to
# row 16 before filling in the values, so the cells reuse the existing style
# rather than Excel minting new cellXfs entries.
$ws.Range("A15:C15").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)

$ws.Range("A16").Value = "Borrowing or Referencing"
$ws.Range("B16").Value = "&"
$ws.Range("C16").Value = "&"

# --- Row 17: "Dereferencing" -------------------------------------------------
# Copy the formatting of row 14 (A/B plain style index 2, C wrap-text style
# index 3) down to row 17.
$ws.Range("A14:C14").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)

# Write B17 ("*") before A17 ("Dereferencing") so the shared-string table
# gets the same new-string insertion order as the source workbook.
$ws.Range("B17").Value = "*"
$ws.Range("A17").Value = "Dereferencing"
$ws.Range("C17").Value = "*`nAutomatically gets dereferenced in case of reference to a reference. No need to write like this **. Instead just use *.`nFollow this file - ""./tuts/ownership/dereference_2.rs"""

$ws.Rows.Item(17).RowHeight = 71

$excel.CutCopyMode = 0

# --- dimension / view ---------------------------------------------------
$ws.Range("A17").Select() | Out-Null

$wb.Save() | Out-Null
